$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = "@"
$r.Value = '67.449.00'
$r.Style = "Normal"
$ws.Range('E2').Value = '  -0.02%  '
$r = $ws.Range('D3')
$r.NumberFormat = "@"
$r.Value = '3.326.54'
$r.Style = "Normal"
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').Value = '  +0.07%  '
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '581.66'
$r.Style = "Normal"
$ws.Range('E5').Value = '  -0.24%  '
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '176.39'
$r.Style = "Normal"
$ws.Range('E6').Value = '  -3.05%  '
$r = $ws.Range('D7')
$r.NumberFormat = "@"
$r.Value = '0.999'
$r.Style = "Normal"
$ws.Range('E7').Value = '  -0.13%  '
$r = $ws.Range('D8')
$r.NumberFormat = "@"
$r.Value = '0.590'
$r.Style = "Normal"
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '3.323.46'
$r.Style = "Normal"
$ws.Range('E9').Value = '  +1.19%  '
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '0.179'
$r.Style = "Normal"
$ws.Range('E10').Value = '  +0.39%  '
$r = $ws.Range('D11')
$r.NumberFormat = "@"
$r.Value = '0.577'
$r.Style = "Normal"
$ws.Range('E11').Value = '  +0.30%  '
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '45.54'
$r.Style = "Normal"
$ws.Range('E12').Value = '  -1.39%  '
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '0.0000270'
$r.Style = "Normal"
$ws.Range('E13').Value = '  -1.77%  '
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '659.70'
$r.Style = "Normal"
$ws.Range('E14').Value = '  +3.79%  '
$r = $ws.Range('D15')
$r.NumberFormat = "@"
$r.Value = '3.865.92'
$r.Style = "Normal"
$ws.Range('E15').Value = '  +1.31%  '
$r = $ws.Range('D16')
$r.NumberFormat = "@"
$r.Value = '8.42'
$r.Style = "Normal"
$ws.Range('E16').Value = '  +0.19%  '
$r = $ws.Range('D17')
$r.NumberFormat = "@"
$r.Value = '67.520.51'
$r.Style = "Normal"
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('E18').Value = '  -0.15%  '
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '3.325.92'
$r.Style = "Normal"
$ws.Range('E19').Value = '  +0.95%  '
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '17.39'
$r.Style = "Normal"
$ws.Range('E20').Value = '  -0.98%  '
$r = $ws.Range('D21')
$r.NumberFormat = "@"
$r.Value = '10.97'
$r.Style = "Normal"
$ws.Range('E21').Value = '  +0.94%  '
$r = $ws.Range('D22')
$r.NumberFormat = "@"
$r.Value = '0.889'
$r.Style = "Normal"
$ws.Range('E22').Value = '  -0.43%  '
$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '5.48'
$r.Style = "Normal"
$ws.Range('E23').Value = '  +9.34%  '
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '17.09'
$r.Style = "Normal"
$ws.Range('E24').Value = '  -3.29%  '
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '99.33'
$r.Style = "Normal"
$ws.Range('E25').Value = '  +1.99%  '
$r = $ws.Range('D26')
$r.NumberFormat = "@"
$r.Value = '3.86'
$r.Style = "Normal"
$ws.Range('E26').Value = '  -2.90%  '
$r = $ws.Range('D27')
$r.NumberFormat = "@"
$r.Value = '2.67'
$r.Style = "Normal"
$ws.Range('E27').Value = '  -4.11%  '
$r = $ws.Range('D28')
$r.NumberFormat = "@"
$r.Value = '9.29'
$r.Style = "Normal"
$ws.Range('E28').Value = '  -2.87%  '
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '33.63'
$r.Style = "Normal"
$ws.Range('E29').Value = '  +3.14%  '
$ws.Range('E30').Value = '  +12.49%  '
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '8.45'
$r.Style = "Normal"
$ws.Range('E31').Value = '  -0.76%  '
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '574.23'
$r.Style = "Normal"
$ws.Range('E32').Value = '  -3.42%  '
$r = $ws.Range('D33')
$r.NumberFormat = "@"
$r.Value = '10.98'
$r.Style = "Normal"
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('E34').Value = '  +0.59%  '
$r = $ws.Range('D35')
$r.NumberFormat = "@"
$r.Value = '1.00'
$r.Style = "Normal"
$ws.Range('E35').Value = '  +0.39%  '
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '3.696.57'
$r.Style = "Normal"
$ws.Range('E36').Value = '  -5.88%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '3.39'
$r.Style = "Normal"
$ws.Range('E37').Value = '  -5.85%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$r = $ws.Range('D38')
$r.NumberFormat = "@"
$r.Value = '56.48'
$r.Style = "Normal"
$ws.Range('E38').Value = '  +1.41%  '
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '34.41'
$r.Style = "Normal"
$ws.Range('E39').Value = '  +4.99%  '
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '0.131'
$r.Style = "Normal"
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('E41').Value = '  -2.67%  '
$ws.Range('E42').Value = '  -4.12%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '0.0₃0672'
$r.Style = "Normal"
$ws.Range('E43').Value = '  -1.83%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$r = $ws.Range('D44')
$r.NumberFormat = "@"
$r.Value = '0.334'
$r.Style = "Normal"
$ws.Range('E44').Value = '  -0.19%  '
$r = $ws.Range('D45')
$r.NumberFormat = "@"
$r.Value = '3.27'
$r.Style = "Normal"
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '2.59'
$r.Style = "Normal"
$ws.Range('E47').Value = '  +2.45%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '0.128'
$r.Style = "Normal"
$ws.Range('E48').Value = '  +0.22%  '
$ws.Range('E49').Value = '  -0.14%  '
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '1.35'
$r.Style = "Normal"
$ws.Range('E50').Value = '  +1.38%  '
$r = $ws.Range('D51')
$r.NumberFormat = "@"
$r.Value = '127.93'
$r.Style = "Normal"
$ws.Range('E51').Value = '  -2.03%  '
